$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing "Poids (en kg)" (column C) values for the 2nd (most recent) weigh-in batch
$ws.Range("C53").Value = 79
$ws.Range("C55").Value = 76.8
$ws.Range("C56").Value = 78.3
$ws.Range("C57").Value = 85.1
$ws.Range("C58").Value = 89.4
$ws.Range("C60").Value = 64.5
$ws.Range("C61").Value = 73.9
$ws.Range("C62").Value = 80.7
$ws.Range("C63").Value = 70.1
$ws.Range("C66").Value = 65.7
$ws.Range("C72").Value = 75.1
$ws.Range("C74").Value = 75.8

# Add a new player "Kamal" with a weight entry on row 79
$ws.Range("A79").Value = "Kamal"
$ws.Range("C79").Value = 90

# Update the visible/active selection to match the new state of the sheet
$ws.Application.ActiveWindow.ScrollRow = 51
$ws.Range("F64").Select()
